# "Add files via upload" — add a new tracked task column (J) to the
# "Definition of Done" sheet: header, all rows marked "x", and the
# author credit row, mirroring the existing "MVP 2.0" column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definition of Done")

# New task header in row 1 (same style already present on the blank J1 cell)
$ws.Range("J1").Value = "E-Mail Dozentenumfrage Aufwandbereitschaft"

# Mark every task row complete, matching column I's "x" entries
$ws.Range("J2").Value = "x"
$ws.Range("J3").Value = "x"
$ws.Range("J4").Value = "x"
$ws.Range("J5").Value = "x"
$ws.Range("J6").Value = "x"

# Credit row at the bottom, styled like I8 (centered, same font/border-less look)
$ws.Range("J8").Value = "Heiser/ Netzler"
$ws.Range("J8").HorizontalAlignment = -4108
$ws.Range("J8").IndentLevel = 1

# Cosmetic bits left behind by the resave (selection + row 1 height)
$ws.Rows.Item(1).RowHeight = 179
$ws.Range("I14").Select()
